$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table "Tabela1" currently spans A1:J79 (78 data rows). A new daily
# data row needs to be appended right below the existing data, which
# causes the table (and its AutoFilter) to automatically grow by one row,
# same as if a user typed a new row directly under the table in Excel.
$lo = $ws.ListObjects.Item(1)
$newListRow = $lo.ListRows.Add()
$newRange = $newListRow.Range

# Clone the formatting of the previous data row onto the freshly added
# table row before writing the new values into it.
$lastDataRow = $ws.Range("A79:J79")
$lastDataRow.Copy($newRange)
$excel.CutCopyMode = $false

# Fill in the new day's figures (2020-05-29).
$newRange.Cells.Item(1, 1).Value = 43980
$newRange.Cells.Item(1, 2).Value = 78529
$newRange.Cells.Item(1, 3).Value = 613
$newRange.Cells.Item(1, 4).Value = 1473
$newRange.Cells.Item(1, 5).Value = 0
$newRange.Cells.Item(1, 6).Value = 7
$newRange.Cells.Item(1, 7).Value = 2
$newRange.Cells.Item(1, 8).Value = 0
$newRange.Cells.Item(1, 9).Value = 108
$newRange.Cells.Item(1, 10).Value = 0

$ws.Range("A80:J80").Select()
